$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("indicators")
$lo = $ws.ListObjects.Item(1)

$row1 = $lo.ListRows.Add()
$row2 = $lo.ListRows.Add()
$row3 = $lo.ListRows.Add()

$ws.Range("A100").Value = "HP_numberofpeople"
$ws.Range("B100").Value = "hpop_healthier_dbl_cntd"
$ws.Range("E100").Value = 1
$ws.Range("H100").Value = 1

$ws.Range("B101").Value = "hpop_healthier_plus_dbl_cntd"
$ws.Range("E101").Value = 1
$ws.Range("H101").Value = 1

$ws.Range("B102").Value = "hpop_healthier_minus_dbl_cntd"
$ws.Range("E102").Value = 1
$ws.Range("H102").Value = 1

# Copy formats from row 99 template cells to rows 100-102
$ws.Range("A99").Copy()
$ws.Range("A100:A102").PasteSpecial(-4122)

$ws.Range("A53").Copy()
$ws.Range("B100:B102").PasteSpecial(-4122)

$ws.Range("C99").Copy()
$ws.Range("C100:C102").PasteSpecial(-4122)

$ws.Range("E99").Copy()
$ws.Range("E100:E102").PasteSpecial(-4122)

$ws.Range("H99").Copy()
$ws.Range("H100:H102").PasteSpecial(-4122)

Write-Output $lo.Range.Address()
